$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.670.44'
$ws.Range('E2').Value = '  -2.07%  '
$ws.Range('D3').Value = '1.758.29'
$ws.Range('E3').Value = '  -2.04%  '
$ws.Range('E4').Value = '  +0.38%  '
$ws.Range('D5').Value = '''326.48'
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('D7').Value = '''0.4447'
$ws.Range('E7').Value = '  -1.59%  '
$ws.Range('E8').Value = '  -0.57%  '
$ws.Range('D9').Value = '''45.95'
$ws.Range('E9').Value = '  +2.77%  '
$ws.Range('D10').Value = '''0.07773'
$ws.Range('E10').Value = '  +3.06%  '
$ws.Range('E11').Value = '  -1.73%  '
$ws.Range('E12').Value = '  +0.24%  '
$ws.Range('D13').Value = '''21.76'
$ws.Range('E13').Value = '  -3.58%  '
$ws.Range('D14').Value = '''6.199'
$ws.Range('E14').Value = '  -1.46%  '
$ws.Range('D15').Value = '''7.377'
$ws.Range('E15').Value = '  -2.14%  '
$ws.Range('D16').Value = '1.760.44'
$ws.Range('E16').Value = '  -1.57%  '
$ws.Range('D17').Value = '''91.21'
$ws.Range('E17').Value = '  +12.48%  '
$ws.Range('D18').Value = '''0.00001082'
$ws.Range('E18').Value = '  -0.67%  '
$ws.Range('D19').Value = '''0.06249'
$ws.Range('E19').Value = '  -7.07%  '
$ws.Range('E20').Value = '  +0.28%  '
$ws.Range('D21').Value = '''17.42'
$ws.Range('E21').Value = '  -0.81%  '
$ws.Range('D22').Value = '''6.189'
$ws.Range('E22').Value = '  -2.42%  '
$ws.Range('D23').Value = '''0.5308'
$ws.Range('E23').Value = '  -3.29%  '
$ws.Range('D24').Value = '27.709.74'
$ws.Range('E24').Value = '  -1.89%  '
$ws.Range('D25').Value = '''11.66'
$ws.Range('E25').Value = '  -0.90%  '
$ws.Range('D26').Value = '''2.341'
$ws.Range('E26').Value = '  -3.28%  '
$ws.Range('E27').Value = '  +1.50%  '
$ws.Range('D28').Value = '''153.64'
$ws.Range('E28').Value = '  +1.22%  '
$ws.Range('D29').Value = '''2.358'
$ws.Range('E29').Value = '  +0.24%  '
$ws.Range('D30').Value = '1.959.33'
$ws.Range('E30').Value = '  -1.39%  '
$ws.Range('D31').Value = '''129.32'
$ws.Range('E31').Value = '  -2.71%  '
$ws.Range('E32').Value = '  -1.57%  '
$ws.Range('D33').Value = '''5.778'
$ws.Range('E33').Value = '  -0.69%  '
$ws.Range('D34').Value = '''0.09274'
$ws.Range('E34').Value = '  -1.56%  '
$ws.Range('D35').Value = '''3.693'
$ws.Range('E35').Value = '  -8.33%  '
$ws.Range('D36').Value = '''12.75'
$ws.Range('E36').Value = '  +5.33%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '''0.02341'
$ws.Range('E37').Value = '  +0.62%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').Value = '''0.2191'
$ws.Range('E38').Value = '  -5.69%  '
$ws.Range('D39').Value = '''0.6506'
$ws.Range('E39').Value = '  -0.82%  '
$ws.Range('D40').Value = '''5.094'
$ws.Range('E40').Value = '  -1.41%  '
$ws.Range('D41').Value = '''0.06126'
$ws.Range('E41').Value = '  -3.34%  '
$ws.Range('D42').Value = '''1.192'
$ws.Range('E42').Value = '  -1.22%  '
$ws.Range('D43').Value = '''8.016'
$ws.Range('E43').Value = '  -3.46%  '
$ws.Range('D44').Value = '''1.417'
$ws.Range('E44').Value = '  -3.54%  '
$ws.Range('E45').Value = '  +0.34%  '
$ws.Range('D46').Value = '''13.80'
$ws.Range('E46').Value = '  -2.30%  '
$ws.Range('D47').Value = '''0.6008'
$ws.Range('E47').Value = '  -1.42%  '
$ws.Range('E48').Value = '  -0.87%  '
$ws.Range('D49').Value = '''125.94'
$ws.Range('E49').Value = '  -3.01%  '
$ws.Range('D50').Value = '''2.002'
$ws.Range('E50').Value = '  -1.14%  '
$ws.Range('D51').Value = '''1.146'
$ws.Range('E51').Value = '  -1.26%  '
